# Applies the "Deploying to gh-pages ... LinuxForHealth/alvearie-fhir-ig"
# update to the StructureDefinition-episode-procedure-type-code workbook:
#   - URL / Version / Date / Publisher metadata bumped (ibm.com -> linuxforhealth.org,
#     7.0.0 -> 8.0.0, new publish date, Alvearie Team -> LinuxForHealth Team)
#   - The duplicated ext-1/ele-1 constraint text that had incorrectly been shown
#     on the root "Extension" row is cleared, since it only really belongs on
#     the "Extension.extension" row further down the table.

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet (Property / Value pairs) ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-procedure-type-code"   # URL
$meta.Range("B3").Value = "8.0.0"                                                                                  # Version
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"                                                              # Date
$meta.Range("B9").Value = "LinuxForHealth Team"                                                                    # Publisher

# --- "Elements" sheet (element definition table) ---
$elements = $wb.Worksheets.Item("Elements")

# "Example" column for the Extension.url row used the same URL literal - keep it in sync.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/episode-procedure-type-code"

# Clear the mis-placed ext-1/ele-1 "Constraint(s)" text from the top-level
# "Extension" row (row 2); the correct home for it is the "Extension.extension"
# row (row 4, column AI), which already carries it and is left untouched.
$elements.Range("AI2").Value = ""
